$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" -> "_FV2310", "_new" -> "_FV2404"
$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2310"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2404"
}

# 2. Turn the used range into an Excel Table
$rng = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
